$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "P220 wraptext:" $ws.Range("P220").WrapText()
Write-Host "AC220 wraptext:" $ws.Range("AC220").WrapText()
Write-Host "P224 wraptext:" $ws.Range("P224").WrapText()
